$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1254.3334
$ws.Range("J17").Value = 1254.3334
$ws.Range("L17").Value = 3763.0002
$ws.Range("N17").Value = -4099.0002
$ws.Range("H19").Value = 177.08333
$ws.Range("I19").Value = 155.14285
$ws.Range("J19").Value = 207.8
$ws.Range("K19").Value = 155.14285
$ws.Range("L19").Value = 207.8
$ws.Range("M19").Value = 19.85714999999999
$ws.Range("N19").Value = -557.8
$ws.Range("H20").Value = 36673.668
$ws.Range("I20").Value = 36673.668
$ws.Range("K20").Value = 36673.668
$ws.Range("M20").Value = -36443.668
$ws.Range("H35").Value = 36673.668
$ws.Range("I35").Value = 36673.668
$ws.Range("K35").Value = 36673.668
$ws.Range("M35").Value = -36294.668
$ws.Range("H64").Value = 3416.9714
$ws.Range("I64").Value = 3119.9092
$ws.Range("J64").Value = 3553.125
$ws.Range("K64").Value = 3119.9092
$ws.Range("L64").Value = 3553.125
$ws.Range("M64").Value = -2871.9092
$ws.Range("N64").Value = -4049.125
$ws.Range("H67").Value = 3416.9714
$ws.Range("I67").Value = 3119.9092
$ws.Range("J67").Value = 3553.125
$ws.Range("K67").Value = 3119.9092
$ws.Range("L67").Value = 3553.125
$ws.Range("M67").Value = -2261.9092
$ws.Range("N67").Value = -5269.125
$ws.Range("H103").Value = 659.8
$ws.Range("I103").Value = 600
$ws.Range("J103").Value = 674.75
$ws.Range("K103").Value = 1800
$ws.Range("L103").Value = 2024.25
$ws.Range("M103").Value = -1214
$ws.Range("N103").Value = -3196.25
$ws.Range("H129").Value = 961.3036
$ws.Range("I129").Value = 554.7143
$ws.Range("K129").Value = 1664.1429
$ws.Range("M129").Value = 3335.8571
$ws.Range("H132").Value = 2398
$ws.Range("I132").Value = 1639
$ws.Range("J132").Value = 3916
$ws.Range("K132").Value = 4917
$ws.Range("L132").Value = 11748
$ws.Range("M132").Value = -2387
$ws.Range("N132").Value = -16808
$ws.Range("H137").Value = 2011.2142
$ws.Range("I137").Value = 2078.0715
$ws.Range("J137").Value = 1944.3572
$ws.Range("K137").Value = 6234.2145
$ws.Range("L137").Value = 5833.071599999999
$ws.Range("M137").Value = -3684.2145
$ws.Range("N137").Value = -10933.0716
$ws.Range("H138").Value = 2634280.2
$ws.Range("J138").Value = 4142.2925
$ws.Range("L138").Value = 12426.8775
$ws.Range("N138").Value = -22706.8775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10101.25
$ws.Range("I32").Value = 9505.745999999999
$ws.Range("K32").Value = 9505.745999999999
$ws.Range("M32").Value = -9218.745999999999
$ws.Range("H45").Value = 1224.75
$ws.Range("J45").Value = 1299.6666
$ws.Range("L45").Value = 1299.6666
$ws.Range("N45").Value = -2053.6666
$ws.Range("H107").Value = 23500
$ws.Range("J107").Value = 23500
$ws.Range("L107").Value = 23500
$ws.Range("N107").Value = -31180
$ws.Range("H110").Value = 1103.6428
$ws.Range("I110").Value = 987.5833
$ws.Range("J110").Value = 1800
$ws.Range("K110").Value = 987.5833
$ws.Range("L110").Value = 1800
$ws.Range("M110").Value = 1057.4167
$ws.Range("N110").Value = -5890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4020.8
$ws.Range("I105").Value = 3672.7273
$ws.Range("J105").Value = 4978
$ws.Range("K105").Value = 3672.7273
$ws.Range("L105").Value = 4978
$ws.Range("M105").Value = -1925.7273
$ws.Range("N105").Value = -8472
$ws.Range("H107").Value = 16552.121
$ws.Range("I107").Value = 20715.385
$ws.Range("J107").Value = 1088.5714
$ws.Range("K107").Value = 20715.385
$ws.Range("L107").Value = 1088.5714
$ws.Range("M107").Value = -18795.385
$ws.Range("N107").Value = -4928.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1876.5883
$ws.Range("I31").Value = 1327.7894
$ws.Range("J31").Value = 3480.7693
$ws.Range("K31").Value = 1327.7894
$ws.Range("L31").Value = 3480.7693
$ws.Range("M31").Value = -1032.7894
$ws.Range("N31").Value = -4070.7693
$ws.Range("H34").Value = 1876.5883
$ws.Range("I34").Value = 1327.7894
$ws.Range("J34").Value = 3480.7693
$ws.Range("K34").Value = 1327.7894
$ws.Range("L34").Value = 3480.7693
$ws.Range("M34").Value = -1125.7894
$ws.Range("N34").Value = -3884.7693
$ws.Range("H94").Value = 2007.421
$ws.Range("I94").Value = 1050
$ws.Range("J94").Value = 2262.7334
$ws.Range("K94").Value = 1050
$ws.Range("L94").Value = 2262.7334
$ws.Range("M94").Value = -599
$ws.Range("N94").Value = -3164.7334
$ws.Range("H99").Value = 3024.4092
$ws.Range("I99").Value = 3026.85
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 3026.85
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1528.85
$ws.Range("N99").Value = -5996
$ws.Range("H103").Value = 20366.666
$ws.Range("I103").Value = 1100
$ws.Range("K103").Value = 1100
$ws.Range("M103").Value = 72
$ws.Range("H126").Value = 3024.4092
$ws.Range("I126").Value = 3026.85
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 9080.549999999999
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6610.549999999999
$ws.Range("N126").Value = -13940
$ws.Range("H129").Value = 47723
$ws.Range("J129").Value = 47723
$ws.Range("L129").Value = 47723
$ws.Range("N129").Value = -57723

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1455.125
$ws.Range("I5").Value = 1726.25
$ws.Range("K5").Value = 5178.75
$ws.Range("M5").Value = -5066.75
$ws.Range("H104").Value = 8000
$ws.Range("J104").Value = 8000
$ws.Range("L104").Value = 24000
$ws.Range("N104").Value = -29242
$ws.Range("H113").Value = 655.2727
$ws.Range("I113").Value = 624
$ws.Range("J113").Value = 681.3333
$ws.Range("K113").Value = 1872
$ws.Range("L113").Value = 2043.9999
$ws.Range("M113").Value = 298
$ws.Range("N113").Value = -6383.9999
$ws.Range("H131").Value = 3334.3394
$ws.Range("J131").Value = 1638.7115
$ws.Range("L131").Value = 4916.1345
$ws.Range("N131").Value = -14996.1345
$ws.Range("H134").Value = 4656.116
$ws.Range("I134").Value = 1426.6666
$ws.Range("J134").Value = 6386.1787
$ws.Range("K134").Value = 4279.9998
$ws.Range("L134").Value = 19158.5361
$ws.Range("M134").Value = 790.0002000000004
$ws.Range("N134").Value = -29298.5361
$ws.Range("H135").Value = 1455.125
$ws.Range("I135").Value = 1726.25
$ws.Range("K135").Value = 15536.25
$ws.Range("M135").Value = -13001.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2946.5454
$ws.Range("I102").Value = 3230.2856
$ws.Range("J102").Value = 2450
$ws.Range("K102").Value = 3230.2856
$ws.Range("L102").Value = 2450
$ws.Range("M102").Value = -1608.2856
$ws.Range("N102").Value = -5694
$ws.Range("H132").Value = 2638.0476
$ws.Range("I132").Value = 1900.2222
$ws.Range("J132").Value = 3191.4167
$ws.Range("K132").Value = 5700.6666
$ws.Range("L132").Value = 9574.250100000001
$ws.Range("M132").Value = -3170.6666
$ws.Range("N132").Value = -14634.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1800.7693
$ws.Range("I82").Value = 1734.8125
$ws.Range("J82").Value = 1906.3
$ws.Range("K82").Value = 1734.8125
$ws.Range("L82").Value = 1906.3
$ws.Range("M82").Value = -1373.8125
$ws.Range("N82").Value = -2628.3
$ws.Range("H85").Value = 1800.7693
$ws.Range("I85").Value = 1734.8125
$ws.Range("J85").Value = 1906.3
$ws.Range("K85").Value = 1734.8125
$ws.Range("L85").Value = 1906.3
$ws.Range("M85").Value = -486.8125
$ws.Range("N85").Value = -4402.3
$ws.Range("H98").Value = 32000
$ws.Range("J98").Value = 32000
$ws.Range("L98").Value = 32000
$ws.Range("N98").Value = -37990
$ws.Range("H132").Value = 8785.429
$ws.Range("I132").Value = 9875
$ws.Range("K132").Value = 29625
$ws.Range("M132").Value = -27095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 967
$ws.Range("I96").Value = 1039.6
$ws.Range("J96").Value = 604
$ws.Range("K96").Value = 1039.6
$ws.Range("L96").Value = 604
$ws.Range("M96").Value = 333.4000000000001
$ws.Range("N96").Value = -3350
$ws.Range("H122").Value = 2132.56
$ws.Range("I122").Value = 2262.8125
$ws.Range("J122").Value = 1901
$ws.Range("K122").Value = 6788.4375
$ws.Range("L122").Value = 5703
$ws.Range("M122").Value = -4338.4375
$ws.Range("N122").Value = -10603
$ws.Range("H126").Value = 6360.6665
$ws.Range("I126").Value = 8272.817999999999
$ws.Range("J126").Value = 1102.25
$ws.Range("K126").Value = 24818.454
$ws.Range("L126").Value = 3306.75
$ws.Range("M126").Value = -22348.454
$ws.Range("N126").Value = -8246.75
